$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization in the Bivariate Inference Assignment link (week 8 eval, H9)
$ws.Range("H9").Value = "[Bivariate Inference Assignment](hw/Bivariate_Inference.html) (Due 10/25)`n"

# Turn the Moderation Assignment text into a link and update the due date (week 10 eval, H11)
$ws.Range("H11").Value = "[Moderation Assignment](hw/Moderation.html) (Due 11/1)"

# Update the active selection to H11, matching the cell last edited
$ws.Range("H11").Select()
